# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 120 (pushing all subsequent
# rows down by one) on the single data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(120).Insert()

$ws.Cells.Item(120, 1).Value = 11
$ws.Cells.Item(120, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(120, 3).Value = "Bíobío"
$ws.Cells.Item(120, 4).Value = 45001
$ws.Cells.Item(120, 5).Value = 8
$ws.Cells.Item(120, 6).Value = 100112043
$ws.Cells.Item(120, 7).Value = "Pepino ensalada"
$ws.Cells.Item(120, 8).Value = "Sin especificar"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 220
$ws.Cells.Item(120, 11).Value = 9000
$ws.Cells.Item(120, 12).Value = 9500
$ws.Cells.Item(120, 13).Value = 9227
$ws.Cells.Item(120, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(120, 15).Value = "Región Metropolitana"
$ws.Cells.Item(120, 16).Value = 154
$ws.Cells.Item(120, 17).Value = 60
$ws.Cells.Item(120, 18).Value = "Hortaliza"
